$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 152 (shifts existing rows 152-258 down to 153-259,
# mirroring Excel's native Rows.Insert behaviour).
$ws.Rows.Item(152).Insert()

# Populate the newly inserted (blank) row 152 with the new data record.
$ws.Cells.Item(152, 1).Value = 5
$ws.Cells.Item(152, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(152, 3).Value = "Maule"
$ws.Cells.Item(152, 4).Value = 44957
$ws.Cells.Item(152, 5).Value = 7
$ws.Cells.Item(152, 6).Value = 100112017
$ws.Cells.Item(152, 7).Value = "Apio"
$ws.Cells.Item(152, 8).Value = "Americana (o)"
$ws.Cells.Item(152, 9).Value = "Primera"
$ws.Cells.Item(152, 10).Value = 500
$ws.Cells.Item(152, 11).Value = 12000
$ws.Cells.Item(152, 12).Value = 12000
$ws.Cells.Item(152, 13).Value = 12000
$ws.Cells.Item(152, 14).Value = "$/docena de matas"
$ws.Cells.Item(152, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(152, 16).Value = 2000
$ws.Cells.Item(152, 17).Value = 6
$ws.Cells.Item(152, 18).Value = "Hortaliza"
